# fix($q): improve custom components
#
# The "Creating and configuring a Component" slide (slide 4) had the
# file-name labels on its two callout rectangles swapped: the box drawn
# next to the *HTML* markup was labelled "heroDetail.html" while the box
# drawn next to the *JS* controller code was labelled "heroDetail.js" -
# i.e. both were actually pointing at the wrong snippet. Swap the two
# text labels so each rectangle names the file it is really annotating.

$p = $ppt.ActivePresentation

# --- Slide 4: "Rectangle 12" <-> "Rectangle 13" label swap -----------------
$s = $p.Slides.Item(4)

$rect12 = $s.Shapes.Item("Rectangle 12")   # previously "heroDetail.html"
$rect13 = $s.Shapes.Item("Rectangle 13")   # previously "heroDetail.js"

$rect12.TextFrame.TextRange.Text = "heroDetail.js"
$rect13.TextFrame.TextRange.Text = "heroDetail.html"

# --- Refresh the cached "last saved" date/time fields -----------------------
# The handout master and notes master each carry a cached
# datetimeFigureOut field (today's date at the time of the last save).
# Bring both up to date; harmless if the host treats these masters as
# read-only.
try {
    $hm = $p.HandoutMaster
    $hm.Shapes.Item("Espace réservé de la date 2").TextFrame.TextRange.Text = "09/06/2017"
} catch {
}

try {
    $nm = $p.NotesMaster
    $nm.Shapes.Item("Espace réservé de la date 2").TextFrame.TextRange.Text = "09/06/2017"
} catch {
}
